# Weekly refresh of the "Bruselas (repollito)" price series for
# Vega Central Mapocho de Santiago.
#
# Three new weekly observations were inserted into the series (the data is
# ordered with the most-recent entries mixed in rather than strictly sorted
# by date), which pushes the previously-recorded weeks down:
#   - a new row lands at sheet row 8, pushing the old rows 8-17 down to 9-18
#   - a new row lands at what becomes sheet row 19 (after the first shift),
#     pushing the old row 18 down to 19 as well
#   - a brand-new final week is appended as row 22
#
# Net effect: rows 8-19 end up holding the data that used to live one row
# above them (row 8 itself becomes new data), and three wholly new rows
# (20, 21, 22) are appended at the bottom. This script writes the resulting
# cell values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: new weekly entry ------------------------------------------------
$ws.Cells.Item(8, 4).Value  = 44428   # Fecha
$ws.Cells.Item(8, 10).Value = 16      # Volumen
$ws.Cells.Item(8, 11).Value = 25000   # Precio minimo
$ws.Cells.Item(8, 12).Value = 26000   # Precio maximo
$ws.Cells.Item(8, 13).Value = 25500   # Precio promedio ponderado
$ws.Cells.Item(8, 16).Value = 1700    # Precio $/Kg

# --- Row 9 (was row 8) ------------------------------------------------------
$ws.Cells.Item(9, 4).Value  = 44349
$ws.Cells.Item(9, 10).Value = 21
$ws.Cells.Item(9, 11).Value = 24000
$ws.Cells.Item(9, 12).Value = 25000
$ws.Cells.Item(9, 13).Value = 24524
$ws.Cells.Item(9, 16).Value = 1635

# --- Row 10 (was row 9) -----------------------------------------------------
$ws.Cells.Item(10, 9).Value  = "Primera"   # Calidad
$ws.Cells.Item(10, 10).Value = 25
$ws.Cells.Item(10, 11).Value = 14000
$ws.Cells.Item(10, 12).Value = 15000
$ws.Cells.Item(10, 13).Value = 14480
$ws.Cells.Item(10, 16).Value = 965

# --- Row 11 (was row 10) ----------------------------------------------------
$ws.Cells.Item(11, 4).Value  = 44385
$ws.Cells.Item(11, 9).Value  = "Segunda"
$ws.Cells.Item(11, 10).Value = 16
$ws.Cells.Item(11, 11).Value = 12000
$ws.Cells.Item(11, 12).Value = 12000
$ws.Cells.Item(11, 13).Value = 12000
$ws.Cells.Item(11, 16).Value = 800

# --- Row 12 (was row 11) ----------------------------------------------------
$ws.Cells.Item(12, 4).Value  = 44413
$ws.Cells.Item(12, 10).Value = 25
$ws.Cells.Item(12, 13).Value = 24480
$ws.Cells.Item(12, 16).Value = 1632

# --- Row 13 (was row 12) ----------------------------------------------------
$ws.Cells.Item(13, 4).Value  = 44421
$ws.Cells.Item(13, 10).Value = 18

# --- Row 14 (was row 13) ----------------------------------------------------
$ws.Cells.Item(14, 4).Value  = 44400
$ws.Cells.Item(14, 10).Value = 16

# --- Row 15 (was row 14) ----------------------------------------------------
$ws.Cells.Item(15, 4).Value  = 44390
$ws.Cells.Item(15, 10).Value = 34
$ws.Cells.Item(15, 11).Value = 24000
$ws.Cells.Item(15, 12).Value = 25000
$ws.Cells.Item(15, 13).Value = 24500
$ws.Cells.Item(15, 16).Value = 1633

# --- Row 16 (was row 15) ----------------------------------------------------
$ws.Cells.Item(16, 4).Value  = 44383
$ws.Cells.Item(16, 10).Value = 25
$ws.Cells.Item(16, 11).Value = 13000
$ws.Cells.Item(16, 12).Value = 14000
$ws.Cells.Item(16, 13).Value = 13480
$ws.Cells.Item(16, 16).Value = 899

# --- Row 17 (was row 16) ----------------------------------------------------
$ws.Cells.Item(17, 4).Value = 44336

# --- Row 18: new weekly entry (was row 17, now re-purposed) -----------------
$ws.Cells.Item(18, 4).Value  = 44435
$ws.Cells.Item(18, 10).Value = 34
$ws.Cells.Item(18, 11).Value = 24000
$ws.Cells.Item(18, 12).Value = 25000
$ws.Cells.Item(18, 13).Value = 24500
$ws.Cells.Item(18, 16).Value = 1633

# --- Row 19 (was row 18) ----------------------------------------------------
$ws.Cells.Item(19, 4).Value  = 44351
$ws.Cells.Item(19, 10).Value = 34
$ws.Cells.Item(19, 11).Value = 24000
$ws.Cells.Item(19, 12).Value = 25000
$ws.Cells.Item(19, 13).Value = 24500
$ws.Cells.Item(19, 16).Value = 1633

# --- Row 20: new row (holds the data that used to be in row 18) ------------
$ws.Cells.Item(20, 1).Value  = 9
$ws.Cells.Item(20, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(20, 3).Value  = "Metropolitana"
$ws.Cells.Item(20, 4).Value  = 44343
$ws.Cells.Item(20, 5).Value  = 13
$ws.Cells.Item(20, 6).Value  = 100112035
$ws.Cells.Item(20, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(20, 8).Value  = "Sin especificar"
$ws.Cells.Item(20, 9).Value  = "Primera"
$ws.Cells.Item(20, 10).Value = 26
$ws.Cells.Item(20, 11).Value = 23000
$ws.Cells.Item(20, 12).Value = 24000
$ws.Cells.Item(20, 13).Value = 23500
$ws.Cells.Item(20, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(20, 15).Value = "Hijuelas"
$ws.Cells.Item(20, 16).Value = 1567
$ws.Cells.Item(20, 17).Value = 15
$ws.Cells.Item(20, 18).Value = "Hortaliza"

# --- Row 21: new row (holds the data that used to be in row 19) ------------
$ws.Cells.Item(21, 1).Value  = 9
$ws.Cells.Item(21, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(21, 3).Value  = "Metropolitana"
$ws.Cells.Item(21, 4).Value  = 44418
$ws.Cells.Item(21, 5).Value  = 13
$ws.Cells.Item(21, 6).Value  = 100112035
$ws.Cells.Item(21, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(21, 8).Value  = "Sin especificar"
$ws.Cells.Item(21, 9).Value  = "Primera"
$ws.Cells.Item(21, 10).Value = 16
$ws.Cells.Item(21, 11).Value = 25000
$ws.Cells.Item(21, 12).Value = 26000
$ws.Cells.Item(21, 13).Value = 25500
$ws.Cells.Item(21, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(21, 15).Value = "Hijuelas"
$ws.Cells.Item(21, 16).Value = 1700
$ws.Cells.Item(21, 17).Value = 15
$ws.Cells.Item(21, 18).Value = "Hortaliza"

# --- Row 22: brand-new final week -------------------------------------------
$ws.Cells.Item(22, 1).Value  = 9
$ws.Cells.Item(22, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(22, 3).Value  = "Metropolitana"
$ws.Cells.Item(22, 4).Value  = 44432
$ws.Cells.Item(22, 5).Value  = 13
$ws.Cells.Item(22, 6).Value  = 100112035
$ws.Cells.Item(22, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(22, 8).Value  = "Sin especificar"
$ws.Cells.Item(22, 9).Value  = "Primera"
$ws.Cells.Item(22, 10).Value = 34
$ws.Cells.Item(22, 11).Value = 24000
$ws.Cells.Item(22, 12).Value = 25000
$ws.Cells.Item(22, 13).Value = 24500
$ws.Cells.Item(22, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(22, 15).Value = "Hijuelas"
$ws.Cells.Item(22, 16).Value = 1633
$ws.Cells.Item(22, 17).Value = 15
$ws.Cells.Item(22, 18).Value = "Hortaliza"

# The "Fecha" column uses a date-typed number format; make sure the three
# brand-new rows pick it up (existing rows already carry the style).
$ws.Cells.Item(20, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(21, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(22, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
